$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 839, pushing existing rows 839:880 down to 840:881.
$ws.Rows.Item(839).Insert()

# Populate the newly inserted row with the new daily entry. Column A holds a
# date-like label that must stay plain text (as every other row in the
# column does), so force text formatting before assigning it, then drop the
# formatting again so the cell ends up with no explicit style - matching the
# rest of the sheet.
$ws.Cells.Item(839, 1).NumberFormat = "@"
$ws.Cells.Item(839, 1).Value = "2026/02/19"
$ws.Cells.Item(839, 1).ClearFormats()

$ws.Cells.Item(839, 2).Value = "木"
$ws.Cells.Item(839, 3).Value = 0
$ws.Cells.Item(839, 4).Value = 201
